# Add (Hul, 2012) and (Lijnen, 2007) data to vessel size and vessel density
# (Vessel size (adipose) -> sheet1 / Table1, Vessel density (adipose) -> sheet3 / Table13)

$wb = $excel.ActiveWorkbook

# ---- Vessel size (adipose) ------------------------------------------------
$wsSize = $wb.Worksheets.Item("Vessel size (adipose)")
$loSize = $wsSize.ListObjects.Item(1)

$loSize.ListRows.Add() | Out-Null
$loSize.ListRows.Add() | Out-Null

$wsSize.Range("A7").Value = "Hul et al., 2012"
$wsSize.Range("B7").Value = 59
$wsSize.Range("C7").Value = 5.1
$wsSize.Range("D7").Value = 49
$wsSize.Range("E7").Value = 2.8

$wsSize.Range("A8").Value = "Lijnen et al., 2007"
$wsSize.Range("D8").Value = 108
$wsSize.Range("E8").Value = 7.7

$wsSize.Range("E8").Select() | Out-Null

# ---- Vessel density (adipose) ---------------------------------------------
$wsDensity = $wb.Worksheets.Item("Vessel density (adipose)")
$loDensity = $wsDensity.ListObjects.Item(1)

$loDensity.ListRows.Add() | Out-Null
$loDensity.ListRows.Add() | Out-Null

$wsDensity.Range("A7").Value = "Hul et al., 2012"
$wsDensity.Range("B7").Value = 740
$wsDensity.Range("C7").Value = 96
$wsDensity.Range("D7").Value = 400
$wsDensity.Range("E7").Value = 55

$wsDensity.Range("A8").Value = "Lijnen et al., 2007"
$wsDensity.Range("D8").Value = 238
$wsDensity.Range("E8").Value = 16

# Vessel density (adipose) becomes the active sheet / tab (moves tabSelected
# off of "CBM (muscle)" and onto this sheet, matching the new activeTab).
$wsDensity.Activate()
$wsDensity.Range("A8").Select() | Out-Null
